$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Buy Value in GBP"
$ws.Range("G1").Value = "Sell Value in GBP"
$ws.Range("J1").Value = "Fee Value in GBP"

$ws.Columns.Item(4).ColumnWidth = 15.5
$ws.Columns.Item(7).ColumnWidth = 15.33203125
$ws.Columns.Item(10).ColumnWidth = 15.33203125
$ws.Columns.Item(11).ColumnWidth = 14.33203125
$ws.Columns.Item(13).ColumnWidth = 54.1640625
Write-Host "widths set"
